$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 110.8604276666667
$ws.Cells.Item(2, 8).Value = 332.581283
$ws.Cells.Item(2, 9).Value = 0.2509786052589675
$ws.Cells.Item(2, 10).Value = 0.2509786052589675
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 490.031855
$ws.Cells.Item(2, 14).Value = 1470.095565
$ws.Cells.Item(2, 15).Value = 0.6686419015677429
$ws.Cells.Item(2, 16).Value = 0.6686419015677431
$ws.Cells.Item(2, 17).Value = 54325.14101558999
$ws.Cells.Item(2, 18).Value = 488926.2691403099
$ws.Cells.Item(2, 19).Value = 0.167814811873176
$ws.Cells.Item(2, 20).Value = 0.167814811873176

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 110.8604276666667
$ws.Cells.Item(3, 8).Value = 332.581283
$ws.Cells.Item(3, 9).Value = 0.2509786052589675
$ws.Cells.Item(3, 10).Value = 0.2509786052589675
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 62.79306433333334
$ws.Cells.Item(3, 14).Value = 188.379193
$ws.Cells.Item(3, 15).Value = 0.0856802950924601
$ws.Cells.Item(3, 16).Value = 0.08568029509246011
$ws.Cells.Item(3, 17).Value = 6961.265966493847
$ws.Cells.Item(3, 18).Value = 62651.39369844462
$ws.Cells.Item(3, 19).Value = 0.02150392096048239
$ws.Cells.Item(3, 20).Value = 0.02150392096048239

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 110.8604276666667
$ws.Cells.Item(4, 8).Value = 332.581283
$ws.Cells.Item(4, 9).Value = 0.2509786052589675
$ws.Cells.Item(4, 10).Value = 0.2509786052589675
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.5977846666666666
$ws.Cells.Item(4, 14).Value = 1.793354
$ws.Cells.Item(4, 15).Value = 0.0008156691696053909
$ws.Cells.Item(4, 16).Value = 0.000815669169605391
$ws.Cells.Item(4, 17).Value = 66.27066379924244
$ws.Cells.Item(4, 18).Value = 596.4359741931819
$ws.Cells.Item(4, 19).Value = 0.0002047155105403012
$ws.Cells.Item(4, 20).Value = 0.0002047155105403012

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 110.8604276666667
$ws.Cells.Item(5, 8).Value = 332.581283
$ws.Cells.Item(5, 9).Value = 0.2509786052589675
$ws.Cells.Item(5, 10).Value = 0.2509786052589675
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 179.453674
$ws.Cells.Item(5, 14).Value = 538.361022
$ws.Cells.Item(5, 15).Value = 0.2448621341701915
$ws.Cells.Item(5, 16).Value = 0.2448621341701916
$ws.Cells.Item(5, 17).Value = 19894.31104599458
$ws.Cells.Item(5, 18).Value = 179048.7994139512
$ws.Cells.Item(5, 19).Value = 0.06145515691476883
$ws.Cells.Item(5, 20).Value = 0.06145515691476884

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 184.841802
$ws.Cells.Item(6, 8).Value = 554.525406
$ws.Cells.Item(6, 9).Value = 0.4184661617850055
$ws.Cells.Item(6, 10).Value = 0.4184661617850055
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 490.031855
$ws.Cells.Item(6, 14).Value = 1470.095565
$ws.Cells.Item(6, 15).Value = 0.6686419015677429
$ws.Cells.Item(6, 16).Value = 0.6686419015677431
$ws.Cells.Item(6, 17).Value = 90578.37111560271
$ws.Cells.Item(6, 18).Value = 815205.3400404244
$ws.Cells.Item(6, 19).Value = 0.2798040101576808
$ws.Cells.Item(6, 20).Value = 0.2798040101576809

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 184.841802
$ws.Cells.Item(7, 8).Value = 554.525406
$ws.Cells.Item(7, 9).Value = 0.4184661617850055
$ws.Cells.Item(7, 10).Value = 0.4184661617850055
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 62.79306433333334
$ws.Cells.Item(7, 14).Value = 188.379193
$ws.Cells.Item(7, 15).Value = 0.0856802950924601
$ws.Cells.Item(7, 16).Value = 0.08568029509246011
$ws.Cells.Item(7, 17).Value = 11606.78316447526
$ws.Cells.Item(7, 18).Value = 104461.0484802774
$ws.Cells.Item(7, 19).Value = 0.03585430422794841
$ws.Cells.Item(7, 20).Value = 0.03585430422794842

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 184.841802
$ws.Cells.Item(8, 8).Value = 554.525406
$ws.Cells.Item(8, 9).Value = 0.4184661617850055
$ws.Cells.Item(8, 10).Value = 0.4184661617850055
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.5977846666666666
$ws.Cells.Item(8, 14).Value = 1.793354
$ws.Cells.Item(8, 15).Value = 0.0008156691696053909
$ws.Cells.Item(8, 16).Value = 0.000815669169605391
$ws.Cells.Item(8, 17).Value = 110.495594994636
$ws.Cells.Item(8, 18).Value = 994.4603549517238
$ws.Cells.Item(8, 19).Value = 0.0003413299466911306
$ws.Cells.Item(8, 20).Value = 0.0003413299466911306

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 184.841802
$ws.Cells.Item(9, 8).Value = 554.525406
$ws.Cells.Item(9, 9).Value = 0.4184661617850055
$ws.Cells.Item(9, 10).Value = 0.4184661617850055
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 179.453674
$ws.Cells.Item(9, 14).Value = 538.361022
$ws.Cells.Item(9, 15).Value = 0.2448621341701915
$ws.Cells.Item(9, 16).Value = 0.2448621341701916
$ws.Cells.Item(9, 17).Value = 33170.54047768055
$ws.Cells.Item(9, 18).Value = 298534.8642991249
$ws.Cells.Item(9, 19).Value = 0.1024665174526851
$ws.Cells.Item(9, 20).Value = 0.1024665174526851

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 95.23175666666667
$ws.Cells.Item(10, 8).Value = 285.69527
$ws.Cells.Item(10, 9).Value = 0.2155966197102082
$ws.Cells.Item(10, 10).Value = 0.2155966197102082
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 490.031855
$ws.Cells.Item(10, 14).Value = 1470.095565
$ws.Cells.Item(10, 15).Value = 0.6686419015677429
$ws.Cells.Item(10, 16).Value = 0.6686419015677431
$ws.Cells.Item(10, 17).Value = 46666.59437427529
$ws.Cells.Item(10, 18).Value = 419999.3493684776
$ws.Cells.Item(10, 19).Value = 0.1441569337746111
$ws.Cells.Item(10, 20).Value = 0.1441569337746112

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 95.23175666666667
$ws.Cells.Item(11, 8).Value = 285.69527
$ws.Cells.Item(11, 9).Value = 0.2155966197102082
$ws.Cells.Item(11, 10).Value = 0.2155966197102082
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 62.79306433333334
$ws.Cells.Item(11, 14).Value = 188.379193
$ws.Cells.Item(11, 15).Value = 0.0856802950924601
$ws.Cells.Item(11, 16).Value = 0.08568029509246011
$ws.Cells.Item(11, 17).Value = 5979.893822946347
$ws.Cells.Item(11, 18).Value = 53819.04440651711
$ws.Cells.Item(11, 19).Value = 0.01847238199770754
$ws.Cells.Item(11, 20).Value = 0.01847238199770754

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 95.23175666666667
$ws.Cells.Item(12, 8).Value = 285.69527
$ws.Cells.Item(12, 9).Value = 0.2155966197102082
$ws.Cells.Item(12, 10).Value = 0.2155966197102082
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.5977846666666666
$ws.Cells.Item(12, 14).Value = 1.793354
$ws.Cells.Item(12, 15).Value = 0.0008156691696053909
$ws.Cells.Item(12, 16).Value = 0.000815669169605391
$ws.Cells.Item(12, 17).Value = 56.92808391506444
$ws.Cells.Item(12, 18).Value = 512.3527552355799
$ws.Cells.Item(12, 19).Value = 0.0001758555157687548
$ws.Cells.Item(12, 20).Value = 0.0001758555157687548

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 95.23175666666667
$ws.Cells.Item(13, 8).Value = 285.69527
$ws.Cells.Item(13, 9).Value = 0.2155966197102082
$ws.Cells.Item(13, 10).Value = 0.2155966197102082
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 179.453674
$ws.Cells.Item(13, 14).Value = 538.361022
$ws.Cells.Item(13, 15).Value = 0.2448621341701915
$ws.Cells.Item(13, 16).Value = 0.2448621341701916
$ws.Cells.Item(13, 17).Value = 17089.68861530733
$ws.Cells.Item(13, 18).Value = 153807.197537766
$ws.Cells.Item(13, 19).Value = 0.05279144842212076
$ws.Cells.Item(13, 20).Value = 0.05279144842212078

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 50.778675
$ws.Cells.Item(14, 8).Value = 152.336025
$ws.Cells.Item(14, 9).Value = 0.1149586132458188
$ws.Cells.Item(14, 10).Value = 0.1149586132458188
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 490.031855
$ws.Cells.Item(14, 14).Value = 1470.095565
$ws.Cells.Item(14, 15).Value = 0.6686419015677429
$ws.Cells.Item(14, 16).Value = 0.6686419015677431
$ws.Cells.Item(14, 17).Value = 24883.16830469212
$ws.Cells.Item(14, 18).Value = 223948.5147422291
$ws.Cells.Item(14, 19).Value = 0.07686614576227499
$ws.Cells.Item(14, 20).Value = 0.07686614576227502

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 50.778675
$ws.Cells.Item(15, 8).Value = 152.336025
$ws.Cells.Item(15, 9).Value = 0.1149586132458188
$ws.Cells.Item(15, 10).Value = 0.1149586132458188
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 62.79306433333334
$ws.Cells.Item(15, 14).Value = 188.379193
$ws.Cells.Item(15, 15).Value = 0.0856802950924601
$ws.Cells.Item(15, 16).Value = 0.08568029509246011
$ws.Cells.Item(15, 17).Value = 3188.548606036426
$ws.Cells.Item(15, 18).Value = 28696.93745432783
$ws.Cells.Item(15, 19).Value = 0.009849687906321746
$ws.Cells.Item(15, 20).Value = 0.009849687906321748

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 50.778675
$ws.Cells.Item(16, 8).Value = 152.336025
$ws.Cells.Item(16, 9).Value = 0.1149586132458188
$ws.Cells.Item(16, 10).Value = 0.1149586132458188
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.5977846666666666
$ws.Cells.Item(16, 14).Value = 1.793354
$ws.Cells.Item(16, 15).Value = 0.0008156691696053909
$ws.Cells.Item(16, 16).Value = 0.000815669169605391
$ws.Cells.Item(16, 17).Value = 30.35471330865
$ws.Cells.Item(16, 18).Value = 273.19241977785
$ws.Cells.Item(16, 19).Value = 0.0000937681966052042990562
$ws.Cells.Item(16, 20).Value = 0.0000937681966052043261613

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 50.778675
$ws.Cells.Item(17, 8).Value = 152.336025
$ws.Cells.Item(17, 9).Value = 0.1149586132458188
$ws.Cells.Item(17, 10).Value = 0.1149586132458188
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 179.453674
$ws.Cells.Item(17, 14).Value = 538.361022
$ws.Cells.Item(17, 15).Value = 0.2448621341701915
$ws.Cells.Item(17, 16).Value = 0.2448621341701916
$ws.Cells.Item(17, 17).Value = 9112.41978960195
$ws.Cells.Item(17, 18).Value = 82011.77810641757
$ws.Cells.Item(17, 19).Value = 0.02814901138061683
$ws.Cells.Item(17, 20).Value = 0.02814901138061685
